# Update the "Confusion matrix sums / Solution Space Size" table (rows 13-17)
# on the "Projects" sheet with the new Solution Space Size values (SAD to
# Code), and refresh the current selection on both sheets to reflect the
# edited range, per commit "update solution space size SAD to code in
# result files".

$wb = $excel.ActiveWorkbook

$wsProjects = $wb.Worksheets.Item("Projects")
$wsCode     = $wb.Worksheets.Item("SAD-Code")

# --- Update the Solution Space Size values in column C (rows 13, 14, 16) ---
$wsProjects.Range("C13").Value = 46495
$wsProjects.Range("C14").Value = 25727
$wsProjects.Range("C16").Value = 164736

# --- Row heights for rows 13-17 shrink from 14.5 to 13.8, matching the ---
# --- rest of the data rows in this table (rows 4-8).                   ---
$wsProjects.Range("A13:A17").EntireRow.RowHeight = 13.8

# --- Select the corresponding range on the SAD-Code sheet (plus the    ---
# --- existing E4 selection), matching the new selection recorded there.---
# --- (E4 listed first so it remains the active cell, as in the target.)---
$wsCode.Activate()
$wsCode.Range("E4,C13:C17").Select()

# --- Select the updated range on the Projects sheet, matching the new ---
# --- selection recorded for this worksheet in the saved file. This     ---
# --- also restores Projects as the active sheet (tabSelected="true").  ---
$wsProjects.Activate()
$wsProjects.Range("C13:C17").Select()
